$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole used range to Text format so numeric-looking strings
# (e.g. "01", "40", "500", "150") are preserved as text, not converted
# to numbers.
$rng = $ws.Range("A1:E4")
$rng.NumberFormat = "@"

$ws.Range("A1").Value = "類別"
$ws.Range("B1").Value = "編號"
$ws.Range("C1").Value = "品名"
$ws.Range("D1").Value = "單價"
$ws.Range("E1").Value = "單位"

$ws.Range("A2").Value = "蔬菜"
$ws.Range("B2").Value = "01"
$ws.Range("C2").Value = "高麗菜"
$ws.Range("D2").Value = "40"
$ws.Range("E2").Value = "粒"

$ws.Range("A3").Value = "水果"
$ws.Range("B3").Value = "02"
$ws.Range("C3").Value = "草莓"
$ws.Range("D3").Value = "500"
$ws.Range("E3").Value = "箱"

$ws.Range("A4").Value = "堅果"
$ws.Range("B4").Value = "03"
$ws.Range("C4").Value = "杏仁"
$ws.Range("D4").Value = "150"
$ws.Range("E4").Value = "罐"
